$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Date heading
Replace-Text "2024-09-22 Sunday" "2024-09-23 Monday"

# Table answers (simple 1-for-1 text replacements; cell count per row is unchanged)
Replace-Text "13×75=975" "96×41=3936"
Replace-Text "96×85=8160" "53×54=2862"
Replace-Text "59×70=4130" "37×71=2627"
Replace-Text "57×62=3534" "98×87=8526"
Replace-Text "99×29=2871" "46×33=1518"

Replace-Text "81×37=2997" "24×66=1584"
Replace-Text "93×52=4836" "62×77=4774"
Replace-Text "37×19=703" "51×52=2652"
Replace-Text "95×96=9120" "58×85=4930"
Replace-Text "48×72=3456" "63×92=5796"

Replace-Text "29×38=1102" "94×95=8930"
Replace-Text "33×88=2904" "92×66=6072"
Replace-Text "90×45=4050" "27×32=864"
Replace-Text "56×33=1848" "30×84=2520"
Replace-Text "40×28=1120" "87×94=8178"

Replace-Text "14×13=182" "25×95=2375"
Replace-Text "74×90=6660" "82×74=6068"
Replace-Text "14×94=1316" "84×71=5964"
Replace-Text "73×33=2409" "37×17=629"
Replace-Text "52×85=4420" "41×73=2993"

# Last row: first two cells replace 1-for-1; the remaining three cells' text is
# reassigned directly (the diff's insert/delete of <w:tc> blocks around the
# unchanged "87x91=7917" cell is equivalent, at the row level, to relabeling
# cells 3-5 while the cell count stays at 5).
Replace-Text "97×35=3395" "31×35=1085"
Replace-Text "30×60=1800" "99×90=8910"

$t = $d.Tables.Item(1)
$lastRow = $t.Rows.Item($t.Rows.Count)
$lastRow.Cells.Item(3).Range.Text = "11×62=682"
$lastRow.Cells.Item(4).Range.Text = "19×22=418"
$lastRow.Cells.Item(5).Range.Text = "87×91=7917"
